# New crime data collected - weekly CompStat refresh
# Advances the report from Week 21 (5/19/2025-5/25/2025) to Week 22
# (5/26/2025-6/1/2025) and refreshes the crime-complaint figures for
# Rape, Robbery, Fel. Assault, Burglary, Gr. Larceny, G.L.A., TOTAL,
# Petit Larceny, Retail Theft, Misd. Assault and UCR Rape*.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header strings: bump the volume/number and the week-covering dates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# ---------------------------------------------------------------------
# 2) Helper: paste the value+format of a donor cell onto a target cell.
#    Pasting value-only first and then format afterwards is required:
#    this host's PasteSpecial(xlPasteAll) alone will carry the donor's
#    *value* across but not its style, so the two-step dance is needed
#    to land a cell on the shared "placeholder text" style (used for
#    the literal "0" / "***.*" entries) or back on a plain numeric
#    style when a cell stops being blank.
# ---------------------------------------------------------------------
function Copy-ValueAndFormat($dstRow, $dstCol, $srcRow, $srcCol) {
    $src = $ws.Cells.Item($srcRow, $srcCol)
    $dst = $ws.Cells.Item($dstRow, $dstCol)
    $src.Copy()
    $dst.PasteSpecial(-4104)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 3) Cells that flip from a numeric entry to the "0" placeholder text
#    (style copied from C14, a stable untouched "0" placeholder cell)
# ---------------------------------------------------------------------
$zeroPlaceholders = @(
    @{r=15; c=4},
    @{r=17; c=4},
    @{r=18; c=3},
    @{r=27; c=4}
)
foreach ($p in $zeroPlaceholders) {
    Copy-ValueAndFormat $p.r $p.c 14 3
}

# ---------------------------------------------------------------------
# 4) Cells that flip from a numeric entry to the "***.*" placeholder
#    text (style copied from E14, a stable untouched "***.*" cell)
# ---------------------------------------------------------------------
$naPlaceholders = @(
    @{r=15; c=5},
    @{r=17; c=5},
    @{r=27; c=5}
)
foreach ($p in $naPlaceholders) {
    Copy-ValueAndFormat $p.r $p.c 14 5
}

# ---------------------------------------------------------------------
# 5) Cells that flip from the placeholder text back to plain numbers
#    (style copied from J14 / K14, stable untouched numeric cells)
# ---------------------------------------------------------------------
Copy-ValueAndFormat 20 4 14 10
$ws.Cells.Item(20, 4).Value = 1

Copy-ValueAndFormat 20 5 14 11
$ws.Cells.Item(20, 5).Value = 0

# ---------------------------------------------------------------------
# 6) Every other changed cell is a like-for-like numeric update that
#    keeps its existing style, so a plain value assignment suffices.
# ---------------------------------------------------------------------
$numericSimple = @(
    @{r=16; c=6; v=4},
    @{r=16; c=7; v=1},
    @{r=16; c=8; v=300},
    @{r=16; c=9; v=8},
    @{r=16; c=11; v=-20},
    @{r=16; c=12; v=60},
    @{r=16; c=13; v=-27.272727272727},
    @{r=16; c=14; v=-63.636363636363},
    @{r=17; c=6; v=5},
    @{r=17; c=7; v=3},
    @{r=17; c=8; v=66.666666666666},
    @{r=17; c=9; v=42},
    @{r=17; c=11; v=2.439024390243},
    @{r=17; c=12; v=16.666666666666},
    @{r=17; c=13; v=121.052631578947},
    @{r=17; c=14; v=0},
    @{r=18; c=4; v=2},
    @{r=18; c=5; v=-100},
    @{r=18; c=7; v=6},
    @{r=18; c=8; v=-33.333333333333},
    @{r=18; c=10; v=15},
    @{r=18; c=11; v=86.666666666666},
    @{r=18; c=12; v=16.666666666666},
    @{r=18; c=13; v=-44},
    @{r=18; c=14; v=-76.859504132231},
    @{r=19; c=4; v=2},
    @{r=19; c=5; v=100},
    @{r=19; c=6; v=21},
    @{r=19; c=7; v=20},
    @{r=19; c=8; v=5},
    @{r=19; c=9; v=104},
    @{r=19; c=10; v=116},
    @{r=19; c=11; v=-10.344827586206},
    @{r=19; c=12; v=-15.447154471544},
    @{r=19; c=13; v=79.310344827586},
    @{r=19; c=14; v=46.478873239436},
    @{r=20; c=3; v=1},
    @{r=20; c=7; v=3},
    @{r=20; c=8; v=0},
    @{r=20; c=9; v=12},
    @{r=20; c=10; v=17},
    @{r=20; c=11; v=-29.411764705882},
    @{r=20; c=12; v=-64.705882352941},
    @{r=20; c=13; v=-20},
    @{r=20; c=14; v=-95.862068965517},
    @{r=21; c=3; v=7},
    @{r=21; c=4; v=5},
    @{r=21; c=5; v=40},
    @{r=21; c=6; v=37},
    @{r=21; c=7; v=34},
    @{r=21; c=8; v=8.823529411764},
    @{r=21; c=9; v=196},
    @{r=21; c=10; v=203},
    @{r=21; c=11; v=-3.448275862068},
    @{r=21; c=12; v=-12.5},
    @{r=21; c=13; v=26.451612903225},
    @{r=21; c=14; v=-64.298724954462},
    @{r=24; c=3; v=8},
    @{r=24; c=4; v=7},
    @{r=24; c=5; v=14.285714285714},
    @{r=24; c=6; v=41},
    @{r=24; c=7; v=24},
    @{r=24; c=8; v=70.833333333333},
    @{r=24; c=9; v=190},
    @{r=24; c=10; v=169},
    @{r=24; c=11; v=12.426035502958},
    @{r=24; c=12; v=-7.317073170731},
    @{r=24; c=13; v=-7.317073170731},
    @{r=25; c=3; v=8},
    @{r=25; c=4; v=3},
    @{r=25; c=5; v=166.666666666667},
    @{r=25; c=6; v=30},
    @{r=25; c=8; v=400},
    @{r=25; c=9; v=106},
    @{r=25; c=10; v=82},
    @{r=25; c=11; v=29.268292682926},
    @{r=25; c=12; v=2.912621359223},
    @{r=26; c=3; v=2},
    @{r=26; c=4; v=3},
    @{r=26; c=5; v=-33.333333333333},
    @{r=26; c=6; v=15},
    @{r=26; c=7; v=15},
    @{r=26; c=8; v=0},
    @{r=26; c=9; v=91},
    @{r=26; c=10; v=67},
    @{r=26; c=11; v=35.820895522388},
    @{r=26; c=12; v=18.181818181818},
    @{r=26; c=13; v=9.638554216867},
    @{r=28; c=6; v=1},
    @{r=28; c=8; v=0},
    @{r=28; c=12; v=0}
)
foreach ($n in $numericSimple) {
    $ws.Cells.Item($n.r, $n.c).Value = $n.v
}
